# Apply the "Saldo" export update:
#  - EDUARDO's Saldo increases from 150074.25 to 360074.25
#  - A new row for JULIANA (account 004813088, Saldo 38195.84) is inserted
#    right after RENATO's row (row 3) and before THAIS's row (row 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update EDUARDO's Saldo value (row 2, column C)
$ws.Cells.Item(2, 3).Value = 360074.25

# 2. Insert a new row at position 4 (shifts THAIS and everything below down by one)
$ws.Rows.Item(4).Insert()

# 3. Populate the newly inserted row with JULIANA's data.
#    The account number must stay text (leading zeros), so force it with a
#    leading apostrophe the same way a user typing it into Excel would.
$ws.Cells.Item(4, 1).Value = "'004813088"
$ws.Cells.Item(4, 1).ClearFormats()
$ws.Cells.Item(4, 2).Value = "JULIANA"
$ws.Cells.Item(4, 3).Value = 38195.84
